# Word COM-interop script: apply the tracked changes to the document.
#
# 1. Recolour the "Infants aren't allocated a seat..." bullet from the
#    red accent colour (C9211E) to black (000000) - both the paragraph
#    mark run properties and the run itself.
# 2. Give the "Normal" paragraph style explicit suppressAutoHyphens /
#    spacing / justification so it matches the Normal paragraph defaults
#    that used to live only in the document defaults.

$d = $word.ActiveDocument

# --- 1. Recolour the specific paragraph -----------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Infants aren't allocated a seat*") {
        $p.Range.Font.Color = 0   # wdColorBlack / RGB(0,0,0) -> 000000
    }
}

# --- 2. Normal style paragraph formatting ----------------------------
$normal = $d.Styles.Item("Normal")
$normal.ParagraphFormat.Hyphenation = 0        # adds <w:suppressAutoHyphens/>
$normal.ParagraphFormat.SpaceBefore = 0        # adds <w:spacing w:before="0" .../>
$normal.ParagraphFormat.SpaceAfter = 0         # ... w:after="0"/>
$normal.ParagraphFormat.Alignment = 0          # adds <w:jc w:val="left"/>
